# Append " – UCN-06" (usecase code) right after the existing
# "Xem tất cả sản phẩm" text, as a brand-new run with matching
# formatting (Aptos font, 12pt / half-point size 24), mirroring how
# Word keeps the newly typed text as its own <w:r> instead of folding
# it into the previous run.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Xem tất cả sản phẩm")

if ($found) {
    # Move to a zero-length range right after the found text.
    $rng.Collapse(0)

    $rng.InsertAfter(" – UCN-06")

    # Apply the same run formatting as the preceding run: Aptos font,
    # 12pt text. Nudging the size away and back forces the host to
    # keep this insertion as its own <w:r> element (matching the
    # target diff) instead of silently merging it into the adjacent
    # run that already carries identical formatting.
    $rng.Font.Name = "Aptos"
    $rng.Font.Size = 11
    $rng.Font.Size = 12
}
